# RegisterDesign.xlsx update
# - Make Sheet2 the active/selected sheet (was Sheet1)
# - Rework the "List Ops" opcode table on Sheet2:
#     remove the old placeholder Ladd/Lsum/Lalu rows and replace them with a
#     fuller opcode table: Ladd, LLoad, LStore, Land, Lor, Lxor, Lnot
#     (each with its own 8-bit opcode string), plus a leading blank
#     separator row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Rebuild the "List Ops" block (rows 22-27 -> rows 22-29, plus the
# trailing "Jump" block shifts down to rows 31-33).
# ---------------------------------------------------------------------

# Wipe the old block (including formatting) so we can lay it out fresh.
$ws2.Range("A22:R33").Clear()

# Blank separator row under the "List Ops" header (A21).
$ws2.Range("A22").Font.Bold = $true
$r = $ws2.Range("E22")
$r.Value = 0
$ws2.Range("F22").NumberFormat = "@"

# Row 23: Ladd (highlighted header-style row, like the other instruction
# group headers above it).
$ws2.Range("A23:E23").Interior.Color = 65535
$ws2.Range("A23").Value = "Ladd"
$ws2.Range("B23").Value = "Reg w/Mem Address 1"
$ws2.Range("C23").Value = "Reg w/Mem Address 2"
$ws2.Range("E23").Value = 1
$ws2.Range("F23").NumberFormat = "@"

# NOTE: cell values below are intentionally written in this specific
# (not strictly row-major) order so that new entries land in the
# worksheet's shared-string table in the same sequence as the source
# edit (matching the target XML's <si> order exactly).
$ws2.Range("F23").Value = "00010001"

# Row 24/25 names first (LLoad / LStore), then the register-role labels.
$ws2.Range("A24").Value = "LLoad"
$ws2.Range("A25").Value = "LStore"
$ws2.Range("B25").Value = "Source Register"
$ws2.Range("B24").Value = "Destination Register"

$ws2.Range("C24").Value = "Memory Address"
$ws2.Range("E24").Value = 2
$ws2.Range("F24").NumberFormat = "@"
$ws2.Range("F24").Value = "10010010"

$ws2.Range("C25").Value = "Memory Address"
$ws2.Range("E25").Value = 3
$ws2.Range("F25").NumberFormat = "@"

# Row 26: Land
$ws2.Range("A26").Value = "Land"
$ws2.Range("B26").Value = "Source Register"
$ws2.Range("C26").Value = "Source Register"
$ws2.Range("E26").Value = 4
$ws2.Range("F26").NumberFormat = "@"

$ws2.Range("F25").Value = "10010011"
$ws2.Range("F26").Value = "00010100"

# Row 27: Lor
$ws2.Range("A27").Value = "Lor"
$ws2.Range("B27").Value = "Source Register"
$ws2.Range("C27").Value = "Source Register"
$ws2.Range("E27").Value = 5
$ws2.Range("F27").NumberFormat = "@"

# Row 28: Lxor
$ws2.Range("A28").Value = "Lxor"
$ws2.Range("B28").Value = "Source Register"
$ws2.Range("C28").Value = "Source Register"
$ws2.Range("E28").Value = 6
$ws2.Range("F28").NumberFormat = "@"

$ws2.Range("F27").Value = "00010101"
$ws2.Range("F28").Value = "00010110"

# Row 29: Lnot (unary -> only one source register column)
$ws2.Range("A29").Value = "Lnot"
$ws2.Range("B29").Value = "Source Register"
$ws2.Range("E29").Value = 7
$ws2.Range("F29").NumberFormat = "@"
$ws2.Range("F29").Value = "00010111"

# Row 30: blank spacer before the "Jump" section.
$ws2.Range("F30").NumberFormat = "@"

# Row 31: "Jump" section header.
$ws2.Range("A31").Font.Bold = $true
$ws2.Range("A31").Value = "Jump"
$ws2.Range("F31").NumberFormat = "@"

# Row 32: blank spacer.
$ws2.Range("A32").Font.Bold = $true
$ws2.Range("F32").NumberFormat = "@"

# Row 33: jmp (moved down from row 27).
$ws2.Range("A33:E33").Interior.Color = 65535
$ws2.Range("A33").Value = "jmp"
$ws2.Range("B33").Value = "Offset from current position"
$ws2.Range("E33").Value = 0
$ws2.Range("F33").NumberFormat = "@"
$ws2.Range("F33").Value = "00100000"

# ---------------------------------------------------------------------
# Sheet2 becomes the active sheet / tab, with H14 selected.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("H14").Select()
